$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Cup - Cold (12oz)" row (row 3), shifting subsequent rows up
$ws.Rows.Item(3).Delete()

# Append new order line items after the existing data (now ending at row 11)
$newRows = @(
    @("ANPLC4F", "Lid Anchor - 24/32oz (Flat)", "2", "66.89", "133.78"),
    @("ANPLC4LD", "Lid Anchor - 24/32oz (Dome)", "2", "56.82", "113.64"),
    @("ANPM424", "Container - Anchor (24oz)", "2", "47.17", "94.34"),
    @("4541602", "Container - Anchor (16oz)", "2", "43.72", "87.44"),
    @("SLOP325", "Cup - Portion (3.25oz)", "1", "53.50", "53.50"),
    @("TS8", "Tamper Evident - 8oz", "1", "38.29", "38.29"),
    @("TS12", "Tamper Evident - 12oz Square", "2", "38.39", "76.78"),
    @("TS16", "Tamper Evident - 16oz", "2", "41.87", "83.74"),
    @("SAB52032T300", "Lid Salad - 24/32oz Sabert (Round)", "1", "80.81", "80.81"),
    @("SAB12032T300", "Cont Salad - 32oz Sabert (Round)", "1", "88.05", "88.05")
)

$startRow = 12
$endRow = $startRow + $newRows.Length - 1

# Ensure these new cells are stored as text, matching the rest of the sheet
$ws.Range("A$startRow" + ":E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
